$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet restructuring: rename CompositeMilestones -> CompositeAccountMilestones
#    and insert a brand new CompositeMemoMilestones sheet right before "config".
#    NOTE: worksheet object references in this runtime are resolved live by
#    tab position, so any sheet captured *before* an insertion/move that ends
#    up sitting at or after the insertion point will silently point at the
#    wrong sheet afterwards. To avoid that we do ALL the structural
#    operations first, and only re-resolve sheet handles (by name) once the
#    final tab layout is settled, before writing any cell data.
# ---------------------------------------------------------------------------
$wsCompositeAccountTmp = $wb.Worksheets.Item("CompositeMilestones")
$wsCompositeAccountTmp.Name = "CompositeAccountMilestones"

$wsCompositeMemoTmp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsCompositeAccountTmp)
$wsCompositeMemoTmp.Name = "CompositeMemoMilestones"

# Re-resolve every sheet handle we will need, now that the final sheet
# layout (order + names) is in place.
$wsAccountSet = $wb.Worksheets.Item("AccountSet")
$wsBudgetSet = $wb.Worksheets.Item("BudgetSet")
$wsAccountMilestones = $wb.Worksheets.Item("AccountMilestones")
$wsMemoMilestones = $wb.Worksheets.Item("MemoMilestones")
$wsCompositeAccount = $wb.Worksheets.Item("CompositeAccountMilestones")
$wsCompositeMemo = $wb.Worksheets.Item("CompositeMemoMilestones")
$wsConfig = $wb.Worksheets.Item("config")

# ---------------------------------------------------------------------------
# 2. AccountSet: the 5 account rows (2-6) are reversed in order, and the
#    Billing_Start_Dt column (F) switches from text to a real number for the
#    rows that carry interest information.
# ---------------------------------------------------------------------------

# Row 2: test loan: Interest
$wsAccountSet.Range("A2").Value = "test loan: Interest"
$wsAccountSet.Range("B2").Value = 0
$wsAccountSet.Range("C2").Value = 0
$wsAccountSet.Range("D2").Value = 9999
$wsAccountSet.Range("E2").Value = "interest"

# Row 3: test loan: Principal Balance
$wsAccountSet.Range("A3").Value = "test loan: Principal Balance"
$wsAccountSet.Range("B3").Value = 100
$wsAccountSet.Range("C3").Value = 0
$wsAccountSet.Range("D3").Value = 9999
$wsAccountSet.Range("E3").Value = "principal balance"
$wsAccountSet.Range("F3").Value = 20000102
$wsAccountSet.Range("G3").Value = "simple"
$wsAccountSet.Range("H3").Value = 0.01
$wsAccountSet.Range("I3").Value = "daily"
$wsAccountSet.Range("J3").Value = 50

# Row 4: Credit: Prev Stmt Bal (position unchanged, only F4 type flips to numeric)
$wsAccountSet.Range("A4").Value = "Credit: Prev Stmt Bal"
$wsAccountSet.Range("B4").Value = 100
$wsAccountSet.Range("C4").Value = 0
$wsAccountSet.Range("D4").Value = 20000
$wsAccountSet.Range("E4").Value = "prev stmt bal"
$wsAccountSet.Range("F4").Value = 20000102
$wsAccountSet.Range("G4").Value = "compound"
$wsAccountSet.Range("H4").Value = 0.01
$wsAccountSet.Range("I4").Value = "monthly"
$wsAccountSet.Range("J4").Value = 40

# Row 5: Credit: Curr Stmt Bal
$wsAccountSet.Range("A5").Value = "Credit: Curr Stmt Bal"
$wsAccountSet.Range("B5").Value = 100
$wsAccountSet.Range("C5").Value = 0
$wsAccountSet.Range("D5").Value = 20000
$wsAccountSet.Range("E5").Value = "curr stmt bal"
$wsAccountSet.Range("F5").Value = ""
$wsAccountSet.Range("G5").Value = ""
$wsAccountSet.Range("H5").Value = ""
$wsAccountSet.Range("I5").Value = ""
$wsAccountSet.Range("J5").Value = ""

# Row 6: Checking
$wsAccountSet.Range("A6").Value = "Checking"
$wsAccountSet.Range("B6").Value = 2000
$wsAccountSet.Range("C6").Value = 0
$wsAccountSet.Range("D6").Value = 100000
$wsAccountSet.Range("E6").Value = "checking"
$wsAccountSet.Range("F6").Value = ""
$wsAccountSet.Range("G6").Value = ""
$wsAccountSet.Range("H6").Value = ""
$wsAccountSet.Range("I6").Value = ""
$wsAccountSet.Range("J6").Value = ""

# ---------------------------------------------------------------------------
# 3. BudgetSet: Start_Date/End_Date become real numbers (date serials stored
#    as plain numbers) and the memo text for two of the rows now references
#    the new per-account regex milestones.
# ---------------------------------------------------------------------------

$wsBudgetSet.Range("A2").Value = 20000102
$wsBudgetSet.Range("B2").Value = 20000102
$wsBudgetSet.Range("F2").Value = "specific regex"

$wsBudgetSet.Range("A3").Value = 20000102
$wsBudgetSet.Range("B3").Value = 20000102
$wsBudgetSet.Range("F3").Value = "specific regex 2"

$wsBudgetSet.Range("A4").Value = 20000104
$wsBudgetSet.Range("B4").Value = 20000104

# ---------------------------------------------------------------------------
# 4. AccountMilestones: two new test rows.
# ---------------------------------------------------------------------------
$wsAccountMilestones.Range("A2").Value = "test account milestone"
$wsAccountMilestones.Range("B2").Value = "Checking"
$wsAccountMilestones.Range("C2").Value = 0
$wsAccountMilestones.Range("D2").Value = 100

$wsAccountMilestones.Range("A3").Value = "test account milestone"
$wsAccountMilestones.Range("B3").Value = "Checking"
$wsAccountMilestones.Range("C3").Value = 0
$wsAccountMilestones.Range("D3").Value = 200

# ---------------------------------------------------------------------------
# 5. MemoMilestones: two new test rows.
# ---------------------------------------------------------------------------
$wsMemoMilestones.Range("A2").Value = "test memo milestone"
$wsMemoMilestones.Range("B2").Value = "specific regex"

$wsMemoMilestones.Range("A3").Value = "test memo milestone"
$wsMemoMilestones.Range("B3").Value = "specific regex 2"

# ---------------------------------------------------------------------------
# 6. CompositeAccountMilestones (was CompositeMilestones): new headers/shape
#    and two new rows of test data.
# ---------------------------------------------------------------------------
$wsCompositeAccount.Range("A1").Value = "Composite_Milestone_Name"
$wsCompositeAccount.Range("B1").Value = "Account_Name"
$wsCompositeAccount.Range("C1").Value = "Min_Balance"
$wsCompositeAccount.Range("D1").Value = "Max_Balance"
$wsCompositeAccount.Range("E1").Value = ""
$wsCompositeAccount.Range("F1").Value = ""

$wsCompositeAccount.Range("A2").Value = "test composite milestone"
$wsCompositeAccount.Range("B2").Value = "Checking"
$wsCompositeAccount.Range("C2").Value = 0
$wsCompositeAccount.Range("D2").Value = 100

$wsCompositeAccount.Range("A3").Value = "test composite milestone 1"
$wsCompositeAccount.Range("B3").Value = "Checking"
$wsCompositeAccount.Range("C3").Value = 0
$wsCompositeAccount.Range("D3").Value = 100

# ---------------------------------------------------------------------------
# 7. CompositeMemoMilestones: brand new sheet, build headers + 2 data rows.
# ---------------------------------------------------------------------------
$wsCompositeMemo.Range("A1").Value = "Composite_Milestone_Name"
$wsCompositeMemo.Range("A1").Font.Bold = $true
$wsCompositeMemo.Range("A1").HorizontalAlignment = -4108
$wsCompositeMemo.Range("A1").VerticalAlignment = -4160
$wsCompositeMemo.Range("A1").Borders.LineStyle = 1

$wsCompositeMemo.Range("B1").Value = "Milestone_Name"
$wsCompositeMemo.Range("B1").Font.Bold = $true
$wsCompositeMemo.Range("B1").HorizontalAlignment = -4108
$wsCompositeMemo.Range("B1").VerticalAlignment = -4160
$wsCompositeMemo.Range("B1").Borders.LineStyle = 1

$wsCompositeMemo.Range("C1").Value = "Memo_Regex"
$wsCompositeMemo.Range("C1").Font.Bold = $true
$wsCompositeMemo.Range("C1").HorizontalAlignment = -4108
$wsCompositeMemo.Range("C1").VerticalAlignment = -4160
$wsCompositeMemo.Range("C1").Borders.LineStyle = 1

$wsCompositeMemo.Range("D1").Value = "Account_Name"
$wsCompositeMemo.Range("D1").Font.Bold = $true
$wsCompositeMemo.Range("D1").HorizontalAlignment = -4108
$wsCompositeMemo.Range("D1").VerticalAlignment = -4160
$wsCompositeMemo.Range("D1").Borders.LineStyle = 1

$wsCompositeMemo.Range("E1").Value = "Min_Balance"
$wsCompositeMemo.Range("E1").Font.Bold = $true
$wsCompositeMemo.Range("E1").HorizontalAlignment = -4108
$wsCompositeMemo.Range("E1").VerticalAlignment = -4160
$wsCompositeMemo.Range("E1").Borders.LineStyle = 1

$wsCompositeMemo.Range("A2").Value = "test composite milestone"
$wsCompositeMemo.Range("D2").Value = "test memo milestone 2"
$wsCompositeMemo.Range("E2").Value = "other specific regex"

$wsCompositeMemo.Range("A3").Value = "test composite milestone 1"
$wsCompositeMemo.Range("D3").Value = "test memo milestone 2"
$wsCompositeMemo.Range("E3").Value = "other specific regex"

# ---------------------------------------------------------------------------
# 8. config: Start_Date_YYYYMMDD/End_Date_YYYYMMDD become real numbers.
# ---------------------------------------------------------------------------
$wsConfig.Range("A2").Value = 20000101
$wsConfig.Range("B2").Value = 20000103
